$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.73"
$ws.Range("E2").Value = "'-3.38%"
$ws.Range("D3").Value = "'51.02"
$ws.Range("E3").Value = "'4.91%"
$ws.Range("D4").Value = "'5.178"
$ws.Range("E4").Value = "'-1.76%"
$ws.Range("D5").Value = "'0.07776"
$ws.Range("E5").Value = "'-3.95%"
$ws.Range("E6").Value = "'-2.08%"
$ws.Range("E7").Value = "'10.84%"
$ws.Range("D8").Value = "'1.565"
$ws.Range("E8").Value = "'-4.79%"
$ws.Range("D9").Value = "'0.1215"
$ws.Range("E9").Value = "'-6.08%"
$ws.Range("D10").Value = "'0.1987"
$ws.Range("E10").Value = "'2.15%"
$ws.Range("D11").Value = "'0.04800"
$ws.Range("E11").Value = "'3.84%"
$ws.Range("D12").Value = "'0.09514"
$ws.Range("E12").Value = "'0.74%"
$ws.Range("E13").Value = "'-0.56%"
$ws.Range("D14").Value = "'0.001271"
$ws.Range("E14").Value = "'-3.95%"
$ws.Range("D15").Value = "'0.005788"
$ws.Range("E15").Value = "'-0.10%"
$ws.Range("E16").Value = "'2,015.76%"
$ws.Range("D17").Value = "'3.328"
$ws.Range("E17").Value = "'-0.32%"
$ws.Range("D18").Value = "'2.433"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("D19").Value = "'0.3477"
$ws.Range("E19").Value = "'1.72%"
$ws.Range("D20").Value = "'8.057"
$ws.Range("E20").Value = "'-0.48%"
$ws.Range("D21").Value = "'0.1367"
$ws.Range("E21").Value = "'-1.68%"
$ws.Range("E22").Value = "'-0.97%"
$ws.Range("D23").Value = "'0.04163"
$ws.Range("E23").Value = "'-0.28%"
$ws.Range("D24").Value = "'0.001270"
$ws.Range("E24").Value = "'-2.81%"
$ws.Range("D25").Value = "'0.003947"
$ws.Range("E25").Value = "'-7.13%"
$ws.Range("D26").Value = "'0.0001350"
$ws.Range("E26").Value = "'-0.17%"
$ws.Range("D38").Value = "'0.02602"
$ws.Range("E38").Value = "'-4.42%"
$ws.Range("D39").Value = "'0.06055"
$ws.Range("E39").Value = "'6.34%"
$ws.Range("D40").Value = "'0.01100"
$ws.Range("E40").Value = "'74.57%"
$ws.Range("D41").Value = "'0.007940"
$ws.Range("E41").Value = "'2.18%"
$ws.Range("D42").Value = "'0.1422"
$ws.Range("E42").Value = "'-1.39%"
$ws.Range("D43").Value = "'0.008415"
$ws.Range("E43").Value = "'9.43%"
$ws.Range("D44").Value = "'0.008344"
$ws.Range("E44").Value = "'3.00%"
$ws.Range("D45").Value = "'0.3376"
$ws.Range("E45").Value = "'5.72%"
$ws.Range("D46").Value = "'0.00007258"
$ws.Range("E46").Value = "'5.79%"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.05317"
$ws.Range("E48").Value = "'-20.83%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.002619"
$ws.Range("E49").Value = "'-34.56%"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("E51").Value = "'-0.08%"
